# Apply "Add files via upload" edit:
#  - update existing row 218 (D/E/G values; F recalculates via existing shared formula)
#  - append 14 new data rows (219-232) with matching values/format
#  - rebuild the F219:F232 shared ABS() formula
#  - leave the selection on Q225, matching the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Update existing row 218 ----
$ws.Cells.Item(218, 4).Value = 55     # D218
$ws.Cells.Item(218, 5).Value = 75     # E218
$ws.Cells.Item(218, 7).Value = 0.25   # G218
# F218 is an existing shared formula (ABS(D218-E218)); Excel recalculates it automatically.

# ---- 2. Prepare formatting for the new rows by copying row 218's formats down ----
$fmtSrc = $ws.Range("A218:T218")
[void]$fmtSrc.Copy()
$fmtDst = $ws.Range("A219:T232")
[void]$fmtDst.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---- 3. New row data: Date,Plant_Type,Plant_Size,Low,High,Temp_Diff,Rain,Growth,Pruned,Quadrant,Shade,UV,Humidity,Dew_Point,Pressure,Wind_Gust,Cloud_Cover,Visibility,AQI,Pollen ----
$csv = @"
45818,Flowering,Large,55,75,20,0.25,0.3,No,2,Neutral,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Nonflowering,Medium,55,75,20,0.25,0.3,No,3,Bright,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Nonflowering,Small,55,75,20,0.25,0.3,No,3,Neutral,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Nonflowering,Medium,55,75,20,0.25,0.4,No,3,Bright,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Nonflowering,Medium,55,75,20,0.25,0.5,No,3,Bright,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Nonflowering,Large,55,75,20,0.25,0.75,No,4,Neutral,6,0.76,65,29.77,17,0.75,9.9,50,35
45818,Tree,Medium,55,75,20,0.25,2,No,1,Neutral,6,0.76,65,29.77,17,0.75,9.9,50,35
45819,Flowering,Large,55,75,20,0,0,No,2,Neutral,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Nonflowering,Medium,55,75,20,0,0,No,3,Neutral,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Nonflowering,Small,55,75,20,0,0.1,No,3,Dark,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Nonflowering,Medium,55,75,20,0,0.2,No,3,Bright,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Nonflowering,Medium,55,75,20,0,0.2,No,3,Bright,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Nonflowering,Large,55,75,20,0,0.5,No,4,Bright,7,0.5,59,30.13,15,0.05,9.9,63,14
45819,Tree,Medium,55,75,20,0,1.25,No,1,Neutral,7,0.5,59,30.13,15,0.05,9.9,63,14
"@

$lines = $csv -split "`n"
$startRow = 219
$rowOffset = 0
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $fields = $line -split ","
    $r = $startRow + $rowOffset

    $ws.Cells.Item($r, 1).Value = [double]$fields[0]     # A Date
    $ws.Cells.Item($r, 2).Value = $fields[1]              # B Plant_Type
    $ws.Cells.Item($r, 3).Value = $fields[2]              # C Plant_Size
    $ws.Cells.Item($r, 4).Value = [double]$fields[3]      # D Low
    $ws.Cells.Item($r, 5).Value = [double]$fields[4]      # E High
    # F Temp_Diff formula is handled separately below as a shared formula block
    $ws.Cells.Item($r, 7).Value = [double]$fields[6]      # G Rain
    $ws.Cells.Item($r, 8).Value = [double]$fields[7]      # H Growth
    $ws.Cells.Item($r, 9).Value = $fields[8]              # I Pruned
    $ws.Cells.Item($r, 10).Value = [double]$fields[9]     # J Quadrant
    $ws.Cells.Item($r, 11).Value = $fields[10]            # K Shade
    $ws.Cells.Item($r, 12).Value = [double]$fields[11]    # L UV
    $ws.Cells.Item($r, 13).Value = [double]$fields[12]    # M Humidity
    $ws.Cells.Item($r, 14).Value = [double]$fields[13]    # N Dew_Point
    $ws.Cells.Item($r, 15).Value = [double]$fields[14]    # O Pressure
    $ws.Cells.Item($r, 16).Value = [double]$fields[15]    # P Wind_Gust
    $ws.Cells.Item($r, 17).Value = [double]$fields[16]    # Q Cloud_Cover
    $ws.Cells.Item($r, 18).Value = [double]$fields[17]    # R Visibility
    $ws.Cells.Item($r, 19).Value = [double]$fields[18]    # S AQI
    $ws.Cells.Item($r, 20).Value = [double]$fields[19]    # T Pollen

    $rowOffset = $rowOffset + 1
}

# ---- 4. F219:F232 shared formula ----
$ws.Range("F219:F232").Formula = "=ABS(D219-E219)"

# ---- 5. Restore the saved selection/view state (best effort) ----
$win = $excel.ActiveWindow
$win.ScrollRow = 218
$win.ScrollColumn = 2
$win.Left = 468
$win.Top = 1200
$win.Width = 21156
$win.Height = 6876
[void]$ws.Range("Q225").Select()
